$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.749
$ws.Range("D4").Value = -8.014999999999999
$ws.Range("D5").Value = -8.581999999999999
$ws.Range("B6").Value = 7.295999999999999
$ws.Range("B7").Value = 6.596000000000001
$ws.Range("D8").Value = -8.203999999999999
$ws.Range("B16").Value = 6.657999999999999
$ws.Range("D16").Value = -8.393000000000001
$ws.Range("B20").Value = 6.09
$ws.Range("D22").Value = -8.222000000000001
